$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rows 97 and 98 had their match data (columns F:V) swapped ---
# (the "Indice"/date columns A:E are identical between the two rows, so only
# F:V needs to move)
$row97 = $ws.Range("F97:V97").Value2
$row98 = $ws.Range("F98:V98").Value2

$ws.Range("F97:V97").Value = $row98
$ws.Range("F98:V98").Value = $row97

# --- Step 2: append a new row 99 with the Sittard vs Heracles match ---
# Copy formatting (styles) from the row above so the new row matches the
# existing look (bold/boxed index cell, date-formatted date cell, ...).
$ws.Range("A98:V98").Copy()
$ws.Range("A99:V99").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = "netherlands"
$ws.Cells.Item(99, 3).Value = "eredivisie"
$ws.Cells.Item(99, 4).Value = "2023-2024"
$ws.Cells.Item(99, 5).Value = 45240.83333333334
$ws.Cells.Item(99, 6).Value = "Sittard"
$ws.Cells.Item(99, 7).Value = 4
$ws.Cells.Item(99, 8).Value = "Heracles"
$ws.Cells.Item(99, 9).Value = 1
$ws.Cells.Item(99, 10).Value = 1.68
$ws.Cells.Item(99, 11).Value = "05/11/2023 14:42"
$ws.Cells.Item(99, 12).Value = 1.69
$ws.Cells.Item(99, 13).Value = "10/11/2023 19:59"
$ws.Cells.Item(99, 14).Value = 4.31
$ws.Cells.Item(99, 15).Value = "05/11/2023 14:42"
$ws.Cells.Item(99, 16).Value = 4.05
$ws.Cells.Item(99, 17).Value = "10/11/2023 19:59"
$ws.Cells.Item(99, 18).Value = 4.83
$ws.Cells.Item(99, 19).Value = "05/11/2023 14:42"
$ws.Cells.Item(99, 20).Value = 5.24
$ws.Cells.Item(99, 21).Value = "10/11/2023 19:59"
$ws.Cells.Item(99, 22).Value = "https://www.betexplorer.com/football/netherlands/eredivisie/sittard-heracles/lAnqwMHd/"
